$wb = $excel.ActiveWorkbook

$edits = @(
    @("ALC", "H40", 3800),
    @("ALC", "I40", 3800),
    @("ALC", "J40", 0),
    @("ALC", "K40", 3800),
    @("ALC", "L40", 0),
    @("ALC", "M40", -3625),
    @("ALC", "N40", $null),
    @("ALC", "H112", 9787.666999999999),
    @("ALC", "J112", 12369.857),
    @("ALC", "L112", 37109.571),
    @("ALC", "N112", -39325.571),
    @("ALC", "H116", 35726710),
    @("ALC", "I116", 83340000),
    @("ALC", "K116", 83340000),
    @("ALC", "M116", -83336558),
    @("ALC", "H132", 2364.4707),
    @("ALC", "I132", 1901.8462),
    @("ALC", "J132", 3868),
    @("ALC", "K132", 5705.5386),
    @("ALC", "L132", 11604),
    @("ALC", "M132", -3175.5386),
    @("ALC", "N132", -16664),
    @("ARM", "H32", 4088132.5),
    @("ARM", "I32", 4352456.5),
    @("ARM", "K32", 4352456.5),
    @("ARM", "M32", -4352169.5),
    @("ARM", "H45", 2867.6316),
    @("ARM", "I45", 1996.1428),
    @("ARM", "J45", 3376),
    @("ARM", "K45", 1996.1428),
    @("ARM", "L45", 3376),
    @("ARM", "M45", -1619.1428),
    @("ARM", "N45", -4130),
    @("ARM", "H61", 7992.143),
    @("ARM", "I61", 2257.2727),
    @("ARM", "K61", 2257.2727),
    @("ARM", "M61", -2045.2727),
    @("ARM", "H74", 35166.934),
    @("ARM", "I74", 51918.1),
    @("ARM", "J74", 4710.273),
    @("ARM", "K74", 51918.1),
    @("ARM", "L74", 4710.273),
    @("ARM", "M74", -51044.1),
    @("ARM", "N74", -6458.273),
    @("ARM", "H77", 35166.934),
    @("ARM", "I77", 51918.1),
    @("ARM", "J77", 4710.273),
    @("ARM", "K77", 259590.5),
    @("ARM", "L77", 23551.365),
    @("ARM", "M77", -255222.5),
    @("ARM", "N77", -32287.365),
    @("ARM", "H122", 3961.6),
    @("ARM", "I122", 1576.8334),
    @("ARM", "J122", 6162.923),
    @("ARM", "K122", 4730.5002),
    @("ARM", "L122", 18488.769),
    @("ARM", "M122", -2280.5002),
    @("ARM", "N122", -23388.769),
    @("ARM", "H132", 6103.661),
    @("ARM", "I132", 4557.3096),
    @("ARM", "K132", 13671.9288),
    @("ARM", "M132", -11141.9288),
    @("ARM", "H136", 7992.143),
    @("ARM", "I136", 2257.2727),
    @("ARM", "K136", 6771.8181),
    @("ARM", "M136", -4221.8181),
    @("BSM", "H86", 45571.176),
    @("BSM", "I86", 77911.08),
    @("BSM", "K86", 77911.08),
    @("BSM", "M86", -76788.08),
    @("BSM", "H89", 45571.176),
    @("BSM", "I89", 77911.08),
    @("BSM", "K89", 389555.4),
    @("BSM", "M89", -383939.4),
    @("BSM", "H99", 3954094),
    @("BSM", "I99", 1034.6),
    @("BSM", "K99", 1034.6),
    @("BSM", "M99", 463.4000000000001),
    @("BSM", "H105", 6159.1284),
    @("BSM", "I105", 6878.087),
    @("BSM", "K105", 6878.087),
    @("BSM", "M105", -5131.087),
    @("BSM", "H107", 70319976),
    @("BSM", "I107", 80365256),
    @("BSM", "J107", 3000),
    @("BSM", "K107", 80365256),
    @("BSM", "L107", 3000),
    @("BSM", "M107", -80363336),
    @("BSM", "N107", -6840),
    @("BSM", "H134", 4381.216),
    @("BSM", "I134", 1667.4865),
    @("BSM", "J134", 11553.214),
    @("BSM", "K134", 5002.4595),
    @("BSM", "L134", 34659.642),
    @("BSM", "M134", -2467.4595),
    @("BSM", "N134", -39729.642),
    @("CRP", "H16", 8243.429),
    @("CRP", "I16", 0),
    @("CRP", "K16", 0),
    @("CRP", "M16", $null),
    @("CRP", "H31", 6883.25),
    @("CRP", "I31", 1931.2106),
    @("CRP", "J31", 10127.689),
    @("CRP", "K31", 1931.2106),
    @("CRP", "L31", 10127.689),
    @("CRP", "M31", -1636.2106),
    @("CRP", "N31", -10717.689),
    @("CRP", "H34", 6883.25),
    @("CRP", "I34", 1931.2106),
    @("CRP", "J34", 10127.689),
    @("CRP", "K34", 1931.2106),
    @("CRP", "L34", 10127.689),
    @("CRP", "M34", -1729.2106),
    @("CRP", "N34", -10531.689),
    @("CRP", "H99", 7354.875),
    @("CRP", "I99", 6312),
    @("CRP", "J99", 7503.857),
    @("CRP", "K99", 6312),
    @("CRP", "L99", 7503.857),
    @("CRP", "M99", -4814),
    @("CRP", "N99", -10499.857),
    @("CRP", "H113", 8243.429),
    @("CRP", "I113", 0),
    @("CRP", "K113", 0),
    @("CRP", "M113", $null),
    @("CRP", "H122", 2258.2),
    @("CRP", "I122", 1525.4445),
    @("CRP", "K122", 4576.333500000001),
    @("CRP", "M122", -2126.333500000001),
    @("CRP", "H126", 7354.875),
    @("CRP", "I126", 6312),
    @("CRP", "J126", 7503.857),
    @("CRP", "K126", 18936),
    @("CRP", "L126", 22511.571),
    @("CRP", "M126", -16466),
    @("CRP", "N126", -27451.571),
    @("CRP", "H132", 9253.947),
    @("CRP", "I132", 5978.125),
    @("CRP", "K132", 17934.375),
    @("CRP", "M132", -15404.375),
    @("CUL", "H4", 31708094),
    @("CUL", "I4", 47830644),
    @("CUL", "K4", 143491932),
    @("CUL", "M4", -143491820),
    @("CUL", "H107", 881.7222),
    @("CUL", "I107", 575),
    @("CUL", "J107", 943.06665),
    @("CUL", "K107", 1725),
    @("CUL", "L107", 2829.19995),
    @("CUL", "M107", 195),
    @("CUL", "N107", -6669.19995),
    @("CUL", "H115", 1688.4),
    @("CUL", "I115", 1688.4),
    @("CUL", "K115", 5065.200000000001),
    @("CUL", "M115", -3890.200000000001),
    @("GSM", "H49", 13000),
    @("GSM", "J49", 0),
    @("GSM", "L49", 0),
    @("GSM", "N49", $null),
    @("GSM", "H102", 3654.0588),
    @("GSM", "I102", 3487.9644),
    @("GSM", "K102", 3487.9644),
    @("GSM", "M102", -1865.9644),
    @("GSM", "H113", 6777.6577),
    @("GSM", "J113", 9949.789000000001),
    @("GSM", "L113", 9949.789000000001),
    @("GSM", "N113", -14289.789),
    @("LTW", "H2", 40624.875),
    @("LTW", "I2", 55000),
    @("LTW", "J2", 38571.285),
    @("LTW", "K2", 55000),
    @("LTW", "L2", 38571.285),
    @("LTW", "M2", -54888),
    @("LTW", "N2", -38795.285),
    @("LTW", "H61", 4739.636),
    @("LTW", "I61", 3360.516),
    @("LTW", "J61", 8028.3076),
    @("LTW", "K61", 3360.516),
    @("LTW", "L61", 8028.3076),
    @("LTW", "M61", -3158.516),
    @("LTW", "N61", -8432.3076),
    @("LTW", "H113", 4739.636),
    @("LTW", "I113", 3360.516),
    @("LTW", "J113", 8028.3076),
    @("LTW", "K113", 3360.516),
    @("LTW", "L113", 8028.3076),
    @("LTW", "M113", -1190.516),
    @("LTW", "N113", -12368.3076),
    @("WVR", "H68", 43000),
    @("WVR", "J68", 43000),
    @("WVR", "L68", 43000),
    @("WVR", "N68", -44622),
    @("WVR", "H71", 43000),
    @("WVR", "J71", 43000),
    @("WVR", "L71", 129000),
    @("WVR", "N71", -137112),
    @("WVR", "H81", 13388533),
    @("WVR", "I81", 51299.9),
    @("WVR", "J81", 40063000),
    @("WVR", "K81", 102599.8),
    @("WVR", "L81", 80126000),
    @("WVR", "M81", -101538.8),
    @("WVR", "N81", -80128122),
    @("WVR", "H84", 13388533),
    @("WVR", "I84", 51299.9),
    @("WVR", "J84", 40063000),
    @("WVR", "K84", 512999),
    @("WVR", "L84", 400630000),
    @("WVR", "M84", -507695),
    @("WVR", "N84", -400640608),
    @("WVR", "H95", 400000),
    @("WVR", "J95", 400000),
    @("WVR", "L95", 400000),
    @("WVR", "N95", -405492),
    @("WVR", "H96", 1011.5),
    @("WVR", "I96", 1011.5),
    @("WVR", "K96", 1011.5),
    @("WVR", "M96", 361.5),
    @("WVR", "H100", 712.75),
    @("WVR", "J100", 1014.3333),
    @("WVR", "L100", 2028.6666),
    @("WVR", "N100", -3110.6666),
    @("WVR", "H107", 10417614),
    @("WVR", "I107", 440.25),
    @("WVR", "J107", 20834788),
    @("WVR", "K107", 1320.75),
    @("WVR", "L107", 62504364),
    @("WVR", "M107", 599.25),
    @("WVR", "N107", -62508204),
    @("WVR", "H126", 1251.6364),
    @("WVR", "I126", 1317.4),
    @("WVR", "J126", 1196.8334),
    @("WVR", "K126", 3952.2),
    @("WVR", "L126", 3590.5002),
    @("WVR", "M126", -1482.2),
    @("WVR", "N126", -8530.5002)
)

foreach ($e in $edits) {
    $sheetName = $e[0]
    $cellRef = $e[1]
    $val = $e[2]
    $ws = $wb.Worksheets.Item($sheetName)
    if ($val -eq $null) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $val
    }
}

Write-Output "Done applying edits: total=$($edits.Count)"